$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New row 13 ("UC Copy Trip") added to Tabelle1 (A6:I12 -> A6:I13).
#    The shared-string "UC Copy Trip" is written first so it lands before
#    "UC Export" in the shared strings table, matching the source order.
# ---------------------------------------------------------------------------
$lo1 = $ws.ListObjects.Item("Tabelle1")
$lo1.ListRows.Add() | Out-Null

$ws.Range("A13").Value = "UC Copy Trip"
$ws.Range("B13").Value = 4
$ws.Range("F13").Formula = "=SUM(B13:E13)"
$ws.Range("G13").Value = 34.2
$ws.Range("H13").Formula = "=FALSE"

# ---------------------------------------------------------------------------
# 2. Row 12 ("UC Export / Copy Trip" -> "UC Export", hours reduced)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "UC Export"
$ws.Range("B12").Value = 4
$ws.Range("G12").Value = 34.2
$ws.Range("I12").ClearContents()
$ws.Rows.Item(12).AutoFit()

# ---------------------------------------------------------------------------
# 3. Tabelle3 (A23:C26 -> A23:C27): insert a row before the totals row so the
#    new "UC Copy Trip" evaluation row can be added, pushing "Total Semester 2"
#    from row 26 down to row 27.
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()
$lo3 = $ws.ListObjects.Item("Tabelle3")
$lo3.Resize($ws.Range("A23:C27"))

$ws.Range("A26").Formula = "=A13"
$ws.Range("B26").Formula = "=G12"
$ws.Range("C26").Formula = "=0.2558*B26+9.3956"

$ws.Range("B27").Formula = "=SUM(B24:B26)"
$ws.Range("C27").Formula = "=SUM(C24:C26)"

# ---------------------------------------------------------------------------
# 4. Sheet view + selection
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C28").Select()

# ---------------------------------------------------------------------------
# 5. Reposition the Time/FP chart so it sits below the (now taller) tables.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 498.9228346456693
$co.Top = 263.82540629921
$co.Width = 367.8391811023622
$co.Height = 230.2248062992126
